$wb = $excel.ActiveWorkbook

# --- Debts sheet: insert new "active" column at the very left (A) ---
$wsDebts = $wb.Worksheets.Item("Debts")
$wsDebts.Range("A1").EntireColumn.Insert()
$wsDebts.Range("A1").Value = "active"

# --- Fixed Assets sheet: insert new "active" column at the very left (A) ---
$wsFA = $wb.Worksheets.Item("Fixed Assets")
$wsFA.Range("A1").EntireColumn.Insert()
$wsFA.Range("A1").Value = "active"

# --- Update on-screen selections to match the edited sheets ---
$wsDebts.Activate()
$wsDebts.Range("A1:A1048576").Select()

# Fixed Assets becomes the active/visible tab, with B10 selected
$wsFA.Activate()
$wsFA.Range("B10").Select()
